$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet held 7 toll-reimbursement records (rows 2-8). This test keeps
# only a single record on row 2 and replaces it with a new entry - the
# other six rows (3-8) are removed outright.
$ws.Range("A3:J8").EntireRow.Delete()

# Replace row 2's data with the new record.
$ws.Range("A2").Value = 1557598
$ws.Range("B2").Value = 45516
$ws.Range("C2").Value = 20
$ws.Range("D2").Value = "NTN&2311"
$ws.Range("E2").Value = 550
$ws.Range("F2").Value = "ROTA ITINERANTE DE 2 DIAS"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 215874

# Leave the same cell selected as in the saved workbook.
$ws.Range("J2").Select()
